$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.435.72"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "'1.841.24"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'225.35"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'0.558"
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'32.03"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "'0.293"
$ws.Range("E9").Value = "  +4.41%  "
$ws.Range("D10").Value = "'0.0713"
$ws.Range("E10").Value = "  +8.72%  "
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "'2.110.32"
$ws.Range("E12").Value = "  +3.89%  "
$ws.Range("D13").Value = "'1.844.80"
$ws.Range("E13").Value = "  +3.84%  "
$ws.Range("D14").Value = "'10.86"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").Value = "'0.649"
$ws.Range("E15").Value = "  +4.13%  "
$ws.Range("D16").Value = "'34.471.96"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "'4.37"
$ws.Range("E17").Value = "  +4.08%  "
$ws.Range("D18").Value = "'69.83"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").Value = "'251.66"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "'0.0₃0802"
$ws.Range("E20").Value = "  +8.81%  "
$ws.Range("D21").Value = "'11.29"
$ws.Range("E21").Value = "  +9.22%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'4.29"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "'161.69"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").Value = "'16.74"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").Value = "'7.26"
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D30").Value = "'0.0535"
$ws.Range("E30").Value = "  +4.85%  "
$ws.Range("D31").Value = "'3.82"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").Value = "'3.62"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "'1.94"
$ws.Range("E34").Value = "  +5.12%  "
$ws.Range("D35").Value = "'1.459.23"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").Value = "'0.648"
$ws.Range("E36").Value = "  +3.91%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.07"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0194"
$ws.Range("E38").Value = "  +3.59%  "
$ws.Range("D39").Value = "'0.970"
$ws.Range("E39").Value = "  +9.31%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.84"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'82.43"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "'2.37"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").Value = "'2.15"
$ws.Range("E43").Value = "  +5.28%  "
$ws.Range("D44").Value = "'6.15"
$ws.Range("E44").Value = "  +6.32%  "
$ws.Range("D45").Value = "'2.006.24"
$ws.Range("E45").Value = "  +3.90%  "
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "'0.0499"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").Value = "'106.97"
$ws.Range("E48").Value = "  +9.04%  "
$ws.Range("D49").Value = "'12.24"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  +7.44%  "
